$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 39.41
$ws.Range("B3").Value = 53.18
$ws.Range("B4").Value = 35.31
$ws.Range("B5").Value = 59.62
$ws.Range("B6").Value = 40.95
$ws.Range("B7").Value = 47.79
$ws.Range("B8").Value = 56.88
$ws.Range("B9").Value = 48.39
$ws.Range("B10").Value = 54.21
$ws.Range("B11").Value = 56.42
$ws.Range("B12").Value = 63.75
$ws.Range("B13").Value = 47.82
$ws.Range("B14").Value = 71.7
$ws.Range("B15").Value = 70.75
$ws.Range("B16").Value = 67.81
$ws.Range("B17").Value = 68.93000000000001
$ws.Range("B18").Value = 55.42
$ws.Range("B19").Value = 68.79000000000001
$ws.Range("B20").Value = 54.45
$ws.Range("B21").Value = 58.09
$ws.Range("B22").Value = 61.66
$ws.Range("B23").Value = 78.03
$ws.Range("B24").Value = 71.78
$ws.Range("B25").Value = 70.12
$ws.Range("B26").Value = 73.56999999999999
$ws.Range("B27").Value = 72.06999999999999
$ws.Range("B28").Value = 69.75
$ws.Range("B29").Value = 65.98999999999999
$ws.Range("B30").Value = 74.19
$ws.Range("B31").Value = 60.09
